$wb = $excel.ActiveWorkbook
$first = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($first)
$newSheet.Name = "Table 1"

# --- Numeric cells first (before any text-number-format is applied) ---
$newSheet.Range("E2").Value = 30.2
$newSheet.Range("E4").Value = 8.6

# --- Apply Text number format to the data block (matches s=38 / numFmtId 49) ---
$newSheet.Range("B2:F6").NumberFormat = "@"
$newSheet.Range("A6").NumberFormat = "@"

# --- Header row ---
$newSheet.Range("B1").Value = "Ever used"
$newSheet.Range("C1").Value = "HKCS Ever used"
$newSheet.Range("D1").Value = "Current user"
$newSheet.Range("E1").Value = "HKCS Current user"
$newSheet.Range("F1").Value = "Dependence in last year"

# --- Column A labels ---
$newSheet.Range("A2").Value = "Alcohol"
$newSheet.Range("A3").Value = "Marijuana"
$newSheet.Range("A4").Value = "Cigarettes"
$newSheet.Range("A5").Value = "Cocaine"
$newSheet.Range("A6").Value = "Heroin"

# --- Row 2: Alcohol ---
$newSheet.Range("B2").Value = "37.9"
$newSheet.Range("C2").Value = "59.2"
$newSheet.Range("D2").Value = "13.3"
$newSheet.Range("F2").Value = "4.5"

# --- Row 3: Marijuana ---
$newSheet.Range("B3").Value = "16.6"
$newSheet.Range("C3").Value = "38.0"
$newSheet.Range("D3").Value = "6.4"
$newSheet.Range("E3").Value = "21.2"
$newSheet.Range("F3").Value = "2.1"

# --- Row 4: Cigarettes ---
$newSheet.Range("B4").Value = "NA"
$newSheet.Range("C4").Value = "20.0"
$newSheet.Range("D4").Value = "0.4"
$newSheet.Range("F4").Value = "1.2 (All tobacco)"

# --- Row 5: Cocaine ---
$newSheet.Range("B5").Value = "0.9"
$newSheet.Range("C5").Value = "5.6"
$newSheet.Range("D5").Value = "0.0"
$newSheet.Range("E5").Value = "NA"
$newSheet.Range("F5").Value = "0.0"

# --- Row 6: Heroin ---
$newSheet.Range("B6").Value = "0.0"
$newSheet.Range("C6").Value = "2.0"
$newSheet.Range("D6").Value = "0.0"
$newSheet.Range("E6").Value = "NA"
$newSheet.Range("F6").Value = "0.0"

# --- Footnotes ---
$newSheet.Range("A8").Value = "HKCS = Healthy Kids Colorado Survey"
$newSheet.Range("A9").Value = "Maximum age at recruitment = 18, so some loss of oldest high school kids"
$newSheet.Range("A10").Value = "Median family income is `$100,000-`$150,000, Median household income in Colorado is `$62,520 (2016 American Community Survey)"
$newSheet.Range("A11").Value = "62.1% of the parents who responded have a bachelor's degree or higher, 38.7% of those 25 years or older in Colorado have a bachelor's degree or higher (2016 American Community Survey)"

# --- Column widths (bestFit approximations from diff) ---
$newSheet.Columns.Item(1).ColumnWidth = 11.1640625
$newSheet.Columns.Item(3).ColumnWidth = 13.6640625
$newSheet.Columns.Item(4).ColumnWidth = 11.1640625
$newSheet.Columns.Item(5).ColumnWidth = 15.83203125
$newSheet.Columns.Item(6).ColumnWidth = 16.83203125

$newSheet.Range("H16").Select() | Out-Null

Write-Host "done"
